$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a weekly price log for "Frutilla" (strawberry) at the
# Femacal de La Calera market. This edit adds a new week's worth of
# records (3 quality grades: Especial / Primera / Segunda) as new rows
# 193-195, pushing all the previously-existing data rows (old 193-241)
# down by three (to new rows 196-244).

# Insert 3 blank rows at 193 - this shifts rows 193:241 down to 196:244,
# carrying all of their existing data/formatting with them (native Excel
# "insert rows" behaviour), and grows the sheet dimension to A1:T244.
$ws.Rows("193:195").Insert()

# Columns that are constant across every data row in this block.
$colA = 3
$colB = "Femacal de La Calera"
$colC = "Coquimbo"
$colE = 5
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100112025
$colJ = "Frutilla"
$colK = "Sin especificar"
$colQ = "`$/bandeja 7 kilos"
$colR = "Provincia de Melipilla"
$colT = 7

# NOTE: named (-Param value) arguments are not reliably bound by this
# interpreter, so this helper takes plain positional parameters instead.
function Set-FrutillaRow($RowNum, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($RowNum, 1).Value = $colA
    $ws.Cells.Item($RowNum, 2).Value = $colB
    $ws.Cells.Item($RowNum, 3).Value = $colC
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 5).Value = $colE
    $ws.Cells.Item($RowNum, 6).Value = $colF
    $ws.Cells.Item($RowNum, 7).Value = $colG
    $ws.Cells.Item($RowNum, 8).Value = $colH
    $ws.Cells.Item($RowNum, 9).Value = $colI
    $ws.Cells.Item($RowNum, 10).Value = $colJ
    $ws.Cells.Item($RowNum, 11).Value = $colK
    $ws.Cells.Item($RowNum, 12).Value = $Calidad
    $ws.Cells.Item($RowNum, 13).Value = $Volumen
    $ws.Cells.Item($RowNum, 14).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 15).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 16).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value = $colQ
    $ws.Cells.Item($RowNum, 18).Value = $colR
    $ws.Cells.Item($RowNum, 19).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value = $colT
}

Set-FrutillaRow 193 44588 "Especial" 240 6000 6500 6250 893
Set-FrutillaRow 194 44588 "Primera"  110 5000 5000 5000 714
Set-FrutillaRow 195 44588 "Segunda"  160 4000 4000 4000 571

Write-Host "Inserted 3 new Frutilla rows (193-195); sheet now spans to row $($ws.Cells.Item($ws.Rows.Count,1).Row)."
